# Update "想去人数" (wanted-to-go count) figures that changed between the
# previous data pull and the latest one (gh-pages output regenerated at
# commit 456a3b4).
#
# Sheet "展览" (exhibitions) and sheet "全部类型" (all types, a combined
# view of every sheet) both list the same events, so each updated number
# needs to be written in both places.

$wb = $excel.ActiveWorkbook

$wsExhibit = $wb.Worksheets.Item("展览")
$wsAll     = $wb.Worksheets.Item("全部类型")

# Sheet "展览": F2, F4, F6, F9, F11
$wsExhibit.Range("F2").Value  = 2701
$wsExhibit.Range("F4").Value  = 19300
$wsExhibit.Range("F6").Value  = 2184
$wsExhibit.Range("F9").Value  = 417
$wsExhibit.Range("F11").Value = 226

# Sheet "全部类型": F7, F9, F15, F19, F21 (mirrors the rows above)
$wsAll.Range("F7").Value  = 2701
$wsAll.Range("F9").Value  = 19300
$wsAll.Range("F15").Value = 2184
$wsAll.Range("F19").Value = 417
$wsAll.Range("F21").Value = 226
